$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while forcing text storage (no numeric
# auto-conversion) and without leaving a permanent style change on the cell.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

# Row 2
$ws.Range('D2').Value = '40.143.86'
$ws.Range('E2').Value = '  +1.80%  '

# Row 3
$ws.Range('D3').Value = '2.203.28'
$ws.Range('E3').Value = '  +1.86%  '

# Row 4
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
Set-TextValue 'D5' '228.68'
$ws.Range('E5').Value = '  -0.23%  '

# Row 6
Set-TextValue 'D6' '0.631'
$ws.Range('E6').Value = '  +1.41%  '

# Row 7
Set-TextValue 'D7' '64.01'
$ws.Range('E7').Value = '  +1.40%  '

# Row 8
$ws.Range('E8').Value = '  +0.10%  '

# Row 9
Set-TextValue 'D9' '0.398'
$ws.Range('E9').Value = '  +0.58%  '

# Row 10
Set-TextValue 'D10' '0.0863'
$ws.Range('E10').Value = '  +0.13%  '

# Row 11
$ws.Range('E11').Value = '  +0.24%  '

# Row 12
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').Value = '2.531.77'
$ws.Range('E12').Value = '  +1.97%  '

# Row 13
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D13' '15.98'
$ws.Range('E13').Value = '  -0.31%  '

# Row 14
Set-TextValue 'D14' '22.22'
$ws.Range('E14').Value = '  -0.01%  '

# Row 15
Set-TextValue 'D15' '0.820'
$ws.Range('E15').Value = '  +0.58%  '

# Row 16
Set-TextValue 'D16' '5.59'
$ws.Range('E16').Value = '  +0.24%  '

# Row 17
$ws.Range('D17').Value = '2.208.75'
$ws.Range('E17').Value = '  +2.20%  '

# Row 18
$ws.Range('D18').Value = '40.066.11'
$ws.Range('E18').Value = '  +1.61%  '

# Row 19
$ws.Range('D19').Value = '0.0₃0908'
$ws.Range('E19').Value = '  +6.47%  '

# Row 20
Set-TextValue 'D20' '72.52'
$ws.Range('E20').Value = '  +0.41%  '

# Row 21
Set-TextValue 'D21' '6.11'

# Row 22
Set-TextValue 'D22' '233.30'
$ws.Range('E22').Value = '  +1.98%  '

# Row 23
$ws.Range('E23').Value = '  +0.00%  '

# Row 24
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D24' '2.36'
$ws.Range('E24').Value = '  +0.04%  '

# Row 25
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D25' '2.33'
$ws.Range('E25').Value = '  +0.91%  '

# Row 26
$ws.Range('E26').Value = '  +0.94%  '

# Row 27
Set-TextValue 'D27' '172.01'
$ws.Range('E27').Value = '  -0.08%  '

# Row 28
$ws.Range('E28').Value = '  +1.89%  '

# Row 29
$ws.Range('E29').Value = '  +2.90%  '

# Row 30
Set-TextValue 'D30' '20.12'
$ws.Range('E30').Value = '  +2.06%  '

# Row 31
Set-TextValue 'D31' '2.74'
$ws.Range('E31').Value = '  +4.53%  '

# Row 32
$ws.Range('E32').Value = '  +1.31%  '

# Row 33
Set-TextValue 'D33' '4.59'
$ws.Range('E33').Value = '  -1.43%  '

# Row 34
Set-TextValue 'D34' '4.74'
$ws.Range('E34').Value = '  -1.54%  '

# Row 35
$ws.Range('E35').Value = '  -0.06%  '

# Row 36
Set-TextValue 'D36' '0.0625'
$ws.Range('E36').Value = '  +0.34%  '

# Row 37
Set-TextValue 'D37' '3.92'
$ws.Range('E37').Value = '  +9.16%  '

# Row 38
$ws.Range('E38').Value = '  +0.77%  '

# Row 39
Set-TextValue 'D39' '5.03'
$ws.Range('E39').Value = '  +19.71%  '

# Row 40
$ws.Range('E40').Value = '  +0.14%  '

# Row 41
Set-TextValue 'D41' '103.80'
$ws.Range('E41').Value = '  -0.53%  '

# Row 42
$ws.Range('E42').Value = '  -0.65%  '

# Row 43
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D43' '1.23'
$ws.Range('E43').Value = '  +2.96%  '

# Row 44
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D44' '17.51'
$ws.Range('E44').Value = '  -2.82%  '

# Row 45
$ws.Range('D45').Value = '1.521.89'
$ws.Range('E45').Value = '  -1.05%  '

# Row 46
Set-TextValue 'D46' '8.30'
$ws.Range('E46').Value = '  +5.09%  '

# Row 47
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D47' '0.0931'
$ws.Range('E47').Value = '  -0.21%  '

# Row 48
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D48' '1.11'
$ws.Range('E48').Value = '  +0.40%  '

# Row 49
$ws.Range('E49').Value = '  -0.39%  '

# Row 50
Set-TextValue 'D50' '0.000197'
$ws.Range('E50').Value = '  +33.96%  '

# Row 51
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue 'D51' '49.96'
$ws.Range('E51').Value = '  +7.30%  '
